$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from H1 into the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 8

$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 4

$ws.Range("I9").Value = 4
$ws.Range("J9").Value = 5
